$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 18.372180938720703
$ws.Range("C2").Value = 5.712643623352051
$ws.Range("D2").Value = 12.314848899841309
$ws.Range("E2").Value = 35.0
